$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 17 (index 15): "из них признаны безработными" data columns C/D/E
#    change from text shared-strings ("30,1"/"27,9"/"28,9") to literal
#    numbers (30100/27900/28900).
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 30100
$ws.Range("D17").Value = 27900
$ws.Range("E17").Value = 28900

# ---------------------------------------------------------------------------
# 2) Row 3 columns C/D/E currently hold the "НАСЕЛЕНИЕ" placeholder text
#    (duplicated from B3); restore them to the year headers 2017/2018/2019
#    (stored as text, matching the original shared-string typing) before we
#    overwrite row 2 so those strings remain anchored to row 3 once row 2 is
#    repointed to "Москва". A quoted-text formula forces text typing (plain
#    "2017" would otherwise auto-coerce to a number) and is then baked back
#    down to a static value, so no cell styles/number formats are touched.
# ---------------------------------------------------------------------------
$ws.Range("C3").Formula = "=""2017"""
$ws.Range("C3").Value = $ws.Range("C3").Value
$ws.Range("D3").Formula = "=""2018"""
$ws.Range("D3").Value = $ws.Range("D3").Value
$ws.Range("E3").Formula = "=""2019"""
$ws.Range("E3").Value = $ws.Range("E3").Value

# ---------------------------------------------------------------------------
# 3) Row 2: label the city. B2 is new; C2/D2/E2 switch from the year labels
#    to the city name "Москва".
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Москва"
$ws.Range("C2").Value = "Москва"
$ws.Range("D2").Value = "Москва"
$ws.Range("E2").Value = "Москва"

# ---------------------------------------------------------------------------
# 4) Drop the trailing blank counter rows 77-81 (dataset trimmed to the
#    real 74 data rows -> sheet dimension becomes A1:E76).
# ---------------------------------------------------------------------------
$ws.Range("A77:E81").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 5) Update the view state: selection moves to H16 and the scrolled
#    top-left cell resets back to the sheet origin.
# ---------------------------------------------------------------------------
$ws.Range("H16").Select()
